# plotEIC methods for fGroupsSet
# - Insert a new row for the "getEICsForFGroups" method (implemented, ionize, done)
# - Mark plotEIC as "done"
# - Update the active selection to reflect where the edit was made

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19 ("getFeatures"), shifting everything below it down by one
$ws.Rows(19).Insert()

# Populate the newly inserted row 19 with the new set-method entry
$ws.Range("A19").Value = "getEICsForFGroups"
$ws.Range("D19").Value = "X"
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "X"

# plotEIC (now on row 34 after the insert) is now also marked as done
$ws.Range("G34").Value = "X"

# Reflect the edit location in the sheet's active selection
$ws.Range("G35").Select()
